$d = $word.ActiveDocument

# --- Change 1: remove the empty bullet paragraph that sits right before
#     "Selecionar curso atual;" (paragraph 8 out of 24). Deleting the
#     paragraph mark at the end of the empty paragraph merges it away,
#     leaving "Selecionar curso atual;" as the surviving paragraph. ---
$pEmpty = $d.Paragraphs.Item(8)
$d.Range($pEmpty.Range.End - 1, $pEmpty.Range.End).Delete()

# --- Change 2: drop the "Convidar colega;" and "Iniciar chat;" bullets
#     entirely, and drop the "Enviar mensagem;" bullet's text while
#     keeping its trailing bookmark (_GoBack) attached to the paragraph
#     that used to hold "Cadastrar materia pendente nunca feita;". ---

# Paragraphs, after change 1, read (by index):
#  20 Cadastrar materia pendente nunca feita;
#  21 Convidar colega;
#  22 Iniciar chat;
#  23 Enviar mensagem;   <- carries bookmarkStart/bookmarkEnd "_GoBack"

# Delete "Iniciar chat;" and "Convidar colega;" paragraphs outright
# (highest index first so earlier indices stay valid).
$d.Paragraphs.Item(22).Range.Delete()
$d.Paragraphs.Item(21).Range.Delete()

# Paragraphs are now:
#  20 Cadastrar materia pendente nunca feita;
#  21 Enviar mensagem;   <- still carries the bookmark

# Merge paragraph 21 up into paragraph 20 by deleting paragraph 20's
# paragraph mark (this keeps paragraph 21's mark - and its bookmark -
# as the surviving end-of-paragraph).
$p20 = $d.Paragraphs.Item(20)
$d.Range($p20.Range.End - 1, $p20.Range.End).Delete()

# Now the merged paragraph's text is
# "Cadastrar materia pendente nunca feita;Enviar mensagem;" followed by
# the bookmark. Strip the trailing "Enviar mensagem;" text, leaving the
# bookmark immediately after "Cadastrar materia pendente nunca feita;".
$found = $d.Content.Find.Execute("Enviar mensagem;", $false, $false, $false,
                                  $false, $false, $true, 1, $false, "", 2)
